# Apply cryptos list update (prices / 1h volume % changes, and a few row reorderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a Number by Excel
# (single decimal point, e.g. "0.550" or "20.20") must be forced to Text format
# first so trailing zeros / exact string content survive, matching the original
# inline-string cell type used throughout this sheet.
$textCells = @(
    "D5",
    "D6",
    "D8",
    "D11",
    "D16",
    "D19",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D27",
    "D32",
    "D33",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D42",
    "D46",
    "D50",
    "D51",
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '76.121.51'
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("D3").Value = '2.920.99'
$ws.Range("E3").Value = '  +3.64%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '203.76'
$ws.Range("E5").Value = '  +8.84%  '
$ws.Range("D6").Value = '598.31'
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '0.550'
$ws.Range("E8").Value = '  +0.32%  '
$ws.Range("E9").Value = '  +2.57%  '
$ws.Range("D10").Value = '2.918.25'
$ws.Range("D11").Value = '0.433'
$ws.Range("E11").Value = '  +16.64%  '
$ws.Range("E12").Value = '  +0.54%  '
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("D14").Value = '3.454.60'
$ws.Range("E14").Value = '  +3.52%  '
$ws.Range("D15").Value = '75.961.00'
$ws.Range("E15").Value = '  +1.63%  '
$ws.Range("D16").Value = '28.01'
$ws.Range("E16").Value = '  +4.47%  '
$ws.Range("E17").Value = '  +1.39%  '
$ws.Range("D18").Value = '2.918.85'
$ws.Range("E18").Value = '  +3.72%  '
$ws.Range("D19").Value = '12.96'
$ws.Range("E19").Value = '  +5.31%  '
$ws.Range("E20").Value = '  -2.46%  '
$ws.Range("D21").Value = '372.80'
$ws.Range("E21").Value = '  -1.28%  '
$ws.Range("D22").Value = '2.30'
$ws.Range("E22").Value = '  +1.61%  '
$ws.Range("D23").Value = '4.27'
$ws.Range("E23").Value = '  +5.20%  '
$ws.Range("D24").Value = '71.41'
$ws.Range("E24").Value = '  +0.72%  '
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").Value = '3.071.08'
$ws.Range("E26").Value = '  +3.57%  '
$ws.Range("D27").Value = '4.23'
$ws.Range("E27").Value = '  +1.87%  '
$ws.Range("E28").Value = '  -1.60%  '
$ws.Range("E29").Value = '  +3.82%  '
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("E31").Value = '  +1.01%  '
$ws.Range("D32").Value = '499.77'
$ws.Range("E32").Value = '  -3.36%  '
$ws.Range("D33").Value = '7.72'
$ws.Range("E33").Value = '  -0.23%  '
$ws.Range("E34").Value = '  +2.59%  '
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("B36").Value = 'Cronos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D36").Value = '0.114'
$ws.Range("E36").Value = '  +32.55%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '164.98'
$ws.Range("E37").Value = '  +0.85%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").Value = '20.20'
$ws.Range("E38").Value = '  +1.40%  '
$ws.Range("D39").Value = '19.61'
$ws.Range("E39").Value = '  +1.24%  '
$ws.Range("E40").Value = '  -4.96%  '
$ws.Range("E41").Value = '  +7.14%  '
$ws.Range("D42").Value = '182.58'
$ws.Range("E42").Value = '  -2.45%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("E44").Value = '  +0.17%  '
$ws.Range("E45").Value = '  -0.46%  '
$ws.Range("D46").Value = '40.01'
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("E47").Value = '  -1.97%  '
$ws.Range("E48").Value = '  +1.67%  '
$ws.Range("E49").Value = '  -1.07%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '22.44'
$ws.Range("E50").Value = '  +7.30%  '
$ws.Range("B51").Value = 'Filecoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D51").Value = '3.71'
$ws.Range("E51").Value = '  +0.12%  '
